$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename header cells (row 1) to the new "_per" / lower-cased metric names.
#    Column letters keep their position; only the label text changes for the
#    columns that are being converted from raw counts to population-relative
#    percentages.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "deaths_per"
$ws.Range("I1").Value = "cardiovasc_per"
$ws.Range("J1").Value = "diabetes_per"
$ws.Range("M1").Value = "COPD_per"
$ws.Range("N1").Value = "CKD_per"
$ws.Range("O1").Value = "health_exp"
$ws.Range("Q1").Value = "organ_per"
$ws.Range("R1").Value = "asthma_per"
$ws.Range("S1").Value = "cancer_per"

# ---------------------------------------------------------------------------
# 2. Recompute the data rows (2-54).
#    - Deaths, CKD, Cancer -> divide by Population (col C) and * 100
#    - Cardiovasc, Diabetes, COPD, Organ_transplant, Asthma -> / 1000
#    - HDI, Health_exp -> * 100
#    Columns C,D,E,F,G,H,K,P,T stay untouched.
# ---------------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $population = $ws.Cells.Item($r, 3).Value()

    $deaths = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 2).Value = $deaths / $population * 100

    $cardiovasc = $ws.Cells.Item($r, 9).Value()
    $ws.Cells.Item($r, 9).Value = $cardiovasc / 1000

    $diabetes = $ws.Cells.Item($r, 10).Value()
    $ws.Cells.Item($r, 10).Value = $diabetes / 1000

    $hdi = $ws.Cells.Item($r, 12).Value()
    $ws.Cells.Item($r, 12).Value = $hdi * 100

    $copd = $ws.Cells.Item($r, 13).Value()
    $ws.Cells.Item($r, 13).Value = $copd / 1000

    $ckd = $ws.Cells.Item($r, 14).Value()
    $ws.Cells.Item($r, 14).Value = $ckd / $population * 100

    $healthExp = $ws.Cells.Item($r, 15).Value()
    $ws.Cells.Item($r, 15).Value = $healthExp * 100

    $organ = $ws.Cells.Item($r, 17).Value()
    $ws.Cells.Item($r, 17).Value = $organ / 1000

    $asthma = $ws.Cells.Item($r, 18).Value()
    $ws.Cells.Item($r, 18).Value = $asthma / 1000

    $cancer = $ws.Cells.Item($r, 19).Value()
    $ws.Cells.Item($r, 19).Value = $cancer / $population * 100
}

# ---------------------------------------------------------------------------
# 3. Update the active selection to I5 (matches the saved sheetView state).
# ---------------------------------------------------------------------------
$ws.Range("I5").Select()
